# Simulated Wild Card round and logged it
# Update the "R" (road) row target-depth stats on both the OFF and DEF
# sheets to reflect the additional Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 248
$wsOff.Range("C3").Value = 160
$wsOff.Range("D3").Value = 54
$wsOff.Range("E3").Value = 18

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 293
$wsDef.Range("C3").Value = 206
$wsDef.Range("D3").Value = 70
$wsDef.Range("E3").Value = 32
